$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (45189 -> 45190) for every data row (rows 2 through 407) on each
# automatic refresh of the sheet.
$newValue = 45190
$lastRow = 407

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newValue
}
